$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9-12 (extra cluster rows no longer present after the edit)
$ws.Range("A9:B12").EntireRow.Delete()

# Update cluster names (column A) for the remaining data rows
$ws.Range("A2").Value = "Beyond the City New Years Festival 30 to 31 January Melbourne"
$ws.Range("A3").Value = "Confirmed Omicron Sircuit Bar Fitzroy"
$ws.Range("A4").Value = "Confirmed Omicron Variant The Peel Hotel Collingwood"
$ws.Range("A5").Value = "Melbourne Cricket Ground (MCG)"
$ws.Range("A6").Value = "The Royal Children's Hospital Melbourne Emergency Department Parkville"
$ws.Range("A7").Value = "Werribee Mercy Hospital Emergency Department"
$ws.Range("A8").Value = "Western Health Sunshine Hospital Emergency Department St Albans"

# Update active case counts (column B)
$ws.Range("B2").Value = 36
$ws.Range("B3").Value = 19
$ws.Range("B4").Value = 14
$ws.Range("B5").Value = 22
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 13
$ws.Range("B8").Value = 11
